$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 103, pushing existing rows 103-127
# down to 105-129 (mirrors the target diff, which adds two new weekly
# price records ahead of the existing Frutilla series).
$ws.Rows.Item(103).Insert()
$ws.Rows.Item(103).Insert()

# --- New row 103 -----------------------------------------------------
$ws.Cells.Item(103, 1).Value  = 4
$ws.Cells.Item(103, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(103, 3).Value  = "Los Lagos"
$ws.Cells.Item(103, 4).Value  = 44522
$ws.Cells.Item(103, 5).Value  = 10
$ws.Cells.Item(103, 6).Value  = "Fruta"
$ws.Cells.Item(103, 7).Value  = 100101
$ws.Cells.Item(103, 8).Value  = "Berries"
$ws.Cells.Item(103, 9).Value  = 100112025
$ws.Cells.Item(103, 10).Value = "Frutilla"
$ws.Cells.Item(103, 11).Value = "Sin especificar"
$ws.Cells.Item(103, 12).Value = "Primera"
$ws.Cells.Item(103, 13).Value = 400
$ws.Cells.Item(103, 14).Value = 9000
$ws.Cells.Item(103, 15).Value = 10000
$ws.Cells.Item(103, 16).Value = 9500
$ws.Cells.Item(103, 17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(103, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(103, 19).Value = 1357
$ws.Cells.Item(103, 20).Value = 7

# --- New row 104 -----------------------------------------------------
$ws.Cells.Item(104, 1).Value  = 4
$ws.Cells.Item(104, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(104, 3).Value  = "Los Lagos"
$ws.Cells.Item(104, 4).Value  = 44522
$ws.Cells.Item(104, 5).Value  = 10
$ws.Cells.Item(104, 6).Value  = "Fruta"
$ws.Cells.Item(104, 7).Value  = 100101
$ws.Cells.Item(104, 8).Value  = "Berries"
$ws.Cells.Item(104, 9).Value  = 100112025
$ws.Cells.Item(104, 10).Value = "Frutilla"
$ws.Cells.Item(104, 11).Value = "Sin especificar"
$ws.Cells.Item(104, 12).Value = "Primera"
$ws.Cells.Item(104, 13).Value = 600
$ws.Cells.Item(104, 14).Value = 9000
$ws.Cells.Item(104, 15).Value = 9500
$ws.Cells.Item(104, 16).Value = 9250
$ws.Cells.Item(104, 17).Value = "`$/caja 7 kilos"
$ws.Cells.Item(104, 18).Value = "Región de La Araucanía"
$ws.Cells.Item(104, 19).Value = 1321
$ws.Cells.Item(104, 20).Value = 7
